# Update pl_mw.xlsx results table (rows 2-25) with newly computed values
# for the "case with 380 kV" run. Columns A, G, I, K, L are untouched (0s / index).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B","C","D","E","F","H","J","M","N","O")

# Map of row number -> new values (in the same order as $cols)
$data = @{
    2 = @(0.8122020464435877,0.1345957226481289,0.1849994499780649,0.1509367944983744,1.307652727885483,0.07973214163530429,0.1621527824639308,0.3355715691277439,1.26235982779518,3.087754504164849)
    3 = @(0.7328332203226466,0.1174590017972434,0.1833800931508165,0.1506820903074768,1.305795992869385,0.07973214163530429,0.1626666169639392,0.3150397459795542,1.275168404928138,3.087176000199122)
    4 = @(0.6842055932467019,0.1069039566617107,0.1824459978110298,0.150588278072906,1.30543718257028,0.07973214163530429,0.1630602117256466,0.3025399794180998,1.283539258625073,3.088817609988155)
    5 = @(0.664416904313299,0.1025946173464831,0.1820805468698694,0.1505658086636359,1.305487348835094,0.07973214163530429,0.1632402463017719,0.2974733748442873,1.287077808195264,3.089988244028518)
    6 = @(0.6611326965137891,0.1018785741301258,0.1820207839669337,0.1505630302076781,1.305507538229946,0.07973214163530429,0.1632713272993058,0.2966337172371283,1.287673074963148,3.090212913113533)
    7 = @(0.6839386035113932,0.1068458717813314,0.1824410075814882,0.1505879111991817,1.30543706406246,0.07973214163530429,0.1630625602029845,0.3024715391060084,1.283586465022381,3.088831366950188)
    8 = @(0.7848144384463751,0.1286939708748776,0.1844286393234356,0.1508359987519192,1.306850338438409,0.07973214163530429,0.1623137448547993,0.3284701807484822,1.266671156984081,3.087140316565183)
    9 = @(0.9834310707373675,0.1712677543410734,0.1888015522038273,0.1518181341628058,1.315825453781031,0.07973214163530429,0.1614649933582584,0.3802916761949717,1.237516110773022,3.099692961937166)
    10 = @(1.12980988716015,0.2023739132524724,0.1923012054906224,0.1528409420400436,1.326211335182563,0.07973214163530429,0.161219379845619,0.4188671138120412,1.218541556089214,3.118630717820906)
    11 = @(1.196494336769433,0.2164858720073539,0.1939550189968173,0.1533714777817572,1.331761850221284,0.07973214163530429,0.161189773102862,0.4365235089980288,1.210440145692978,3.129364915390653)
    12 = @(1.22175897567115,0.2218239880407111,0.1945901053931465,0.1535817417944436,1.333982585385598,0.07973214163530429,0.1611903729857858,0.4432248629627793,1.207448586409171,3.133735064077342)
    13 = @(1.216317237252781,0.220674590222643,0.194452936785396,0.1535360416653546,1.333499021585141,0.07973214163530429,0.1611897184468845,0.4417809326093973,1.208089479688056,3.132780284921353)
    14 = @(1.198572625987538,0.2169251592593469,0.1940070914765215,0.1533885888330708,1.331942168375676,0.07973214163530429,0.161189585727243,0.437074529320931,1.210192499604076,3.129718327312474)
    15 = @(1.187705158037318,0.2146277646725423,0.1937351454783993,0.1532994882260965,1.331004035156681,0.07973214163530429,0.1611910426748295,0.434193699746757,1.211490593335377,3.127882571649309)
    16 = @(1.125453737494524,0.2014508691615049,0.1921943630472498,0.152807581503982,1.325865224984611,0.07973214163530429,0.1612229674289409,0.4177153788506942,1.219081673563124,3.117971908942053)
    17 = @(1.087288296636984,0.1933572557215371,0.191264926620903,0.1525225104457633,1.322924348012705,0.07973214163530429,0.161263587558544,0.4076339680218624,1.223874371549478,3.112435263477664)
    18 = @(1.06534564640998,0.1886984252590196,0.1907361604539943,0.1523646885632139,1.32131056690146,0.07973214163530429,0.1612946809670532,0.4018456128867243,1.226680908694924,3.109450176644032)
    19 = @(1.057917839450226,0.1871204151584323,0.1905581310492153,0.1523123085999067,1.320777515849528,0.07973214163530429,0.1613065362001294,0.3998875374042612,1.227639724314258,3.108473714587063)
    20 = @(1.091350139063593,0.1942192089358628,0.1913632647209198,0.1525522210822068,1.323229363468158,0.07973214163530429,0.1612584634587506,0.4087060973625896,1.223359015614697,3.113004002756838)
    21 = @(1.203784313085464,0.2180266168311675,0.1941378081049692,0.1534316454069717,1.332396226789413,0.07973214163530429,0.1611893041387376,0.4384565023099256,1.209572721532183,3.130609406952829)
    22 = @(1.277339966275406,0.2335523460442346,0.1960025307827067,0.1540609566875517,1.339080249918709,0.07973214163530429,0.161212950868979,0.4579889198584937,1.20100715054685,3.143895588780936)
    23 = @(1.238075580807845,0.2252691482454736,0.1950026114224102,0.1537200967275858,1.335449420411095,0.07973214163530429,0.1611940300483923,0.4475560753049947,1.20553806929663,3.136641437087491)
    24 = @(1.089513782602921,0.1938295376376971,0.1913187886829633,0.1525387699969585,1.323091226235576,0.07973214163530429,0.1612607559540962,0.4082213640474066,1.223591848697872,3.112746258993297)
    25 = @(0.9296177607382674,0.1597801857502645,0.1875679833431008,0.1514994746219962,1.312732273432502,0.07973214163530429,0.1616282430516662,0.3661837642558297,1.244973677795059,3.094594007562904)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range("$($cols[$i])$row").Value = $values[$i]
    }
}
